# Update "想去人数" (column F) counts across the workbook's sheets.
# Mirrors the upstream data refresh recorded in the commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F3"  = 534
        "F4"  = 646
        "F7"  = 2156
        "F8"  = 897
        "F9"  = 851
        "F10" = 413
        "F15" = 897
        "F17" = 37
        "F18" = 1776
        "F19" = 45
        "F24" = 1458
        "F26" = 529
        "F27" = 354
        "F28" = 621
        "F29" = 421
        "F30" = 2456
        "F31" = 383
        "F32" = 96
        "F34" = 605
        "F35" = 479
        "F36" = 196
        "F37" = 925
        "F38" = 707
        "F39" = 37
        "F40" = 500
        "F41" = 479
    }
    "演出" = @{
        "F14" = 75
        "F17" = 142
        "F21" = 89
        "F22" = 120
        "F23" = 109
    }
    "本地生活" = @{
        "F3" = 2922
        "F6" = 324
    }
    "全部类型" = @{
        "F3"  = 534
        "F7"  = 646
        "F10" = 897
        "F11" = 851
        "F12" = 413
        "F17" = 897
        "F20" = 37
        "F21" = 324
        "F22" = 1776
        "F23" = 45
        "F29" = 75
        "F30" = 1458
        "F33" = 529
        "F34" = 354
        "F35" = 621
        "F36" = 421
        "F37" = 96
        "F39" = 479
        "F40" = 196
        "F41" = 925
        "F43" = 89
        "F44" = 109
        "F46" = 707
        "F47" = 37
        "F48" = 500
        "F49" = 479
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
